# Update voltage magnitude results (vm_pu) for the "case with 380 kV" run.
# Column B holds the slack-bus setpoint (was 1.05 pu, now 1.02 pu); the
# remaining bus columns (C..F, I..M) are recomputed load-flow results, and
# N2 picks up a tiny floating-point precision change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B="1.02"; C="1.040555794334618"; D="1.047096738552104"; E="1.048319230427385"; F="1.057694054655204"; I="1.035998129261775"; J="1.045641764919755"; K="1.049860330828642"; L="1.051079407225533"; M="1.060428323745304"; N="1.005712725503983" }
    3  = @{ B="1.02"; C="1.041556995282184"; D="1.047870163236246"; E="1.049196448301861"; F="1.058624907041211"; I="1.036153347773119"; J="1.046288373796804"; K="1.050445630435419"; L="1.051768477786647"; M="1.061172765775723" }
    4  = @{ B="1.02"; C="1.042205390041193"; D="1.048370878913515"; E="1.049764887251109"; F="1.059228030536945"; I="1.036252507435273"; J="1.046706721272316"; K="1.05082395628745";  L="1.052214533204774"; M="1.061654639138752" }
    5  = @{ B="1.02"; C="1.04247810627018";  D="1.048581440270978"; E="1.050004054122145"; F="1.059481773805032"; I="1.036293887905323"; J="1.046882581263618"; K="1.050982907191633"; L="1.052402097155591"; M="1.061857258181566" }
    6  = @{ B="1.02"; C="1.042523904188342"; D="1.048616797957882"; E="1.050044222698744"; F="1.059524389528642"; I="1.036300817894674"; J="1.046912108124407"; K="1.051009590013652"; L="1.052433592416829"; M="1.061891281093097" }
    7  = @{ B="1.02"; C="1.042209033574745"; D="1.048373692208672"; E="1.04976808224536";  F="1.059231420321397"; I="1.036253061567139"; J="1.046709071175476"; K="1.050826080580742"; L="1.052217039278698"; M="1.061657346390589" }
    8  = @{ B="1.02"; C="1.040894040359946"; D="1.047358066789423"; E="1.048615519708759"; F="1.058008474192384"; I="1.036050850206405"; J="1.04586029957624";  K="1.050058218225935"; L="1.051312243831235"; M="1.060679875627089" }
    9  = @{ B="1.02"; C="1.038581107598058"; D="1.045570440869253"; E="1.046590895187545"; F="1.05585966999995";  I="1.035684769595195"; J="1.044364299665574"; K="1.048702111481478"; L="1.04971930680116";  M="1.058958795487697" }
    10 = @{ B="1.02"; C="1.037042057483732"; D="1.044380136999011"; E="1.045245480319906"; F="1.054431368238821"; I="1.035434185101928"; J="1.043366780280379"; K="1.04779605912055";  L="1.048658367880479"; M="1.057812380908639" }
    11 = @{ B="1.02"; C="1.036376329879835"; D="1.043865081619201"; E="1.044663944218564"; F="1.053813918111905"; I="1.035324136490078"; J="1.042934810156773"; K="1.047403270820958"; L="1.048199223480448"; M="1.057316214899728" }
    12 = @{ B="1.02"; C="1.03612915352235";  D="1.043673821584395"; E="1.044448092875635"; F="1.053584723232797"; I="1.035283028058852"; J="1.042774352404567"; K="1.047257303352428"; L="1.048028715184055"; M="1.057131953726951" }
    13 = @{ B="1.02"; C="1.036182168949806"; D="1.043714845025951"; E="1.044494386594967"; F="1.053633879341971"; I="1.035291856424111"; J="1.042808771341458"; K="1.047288616950105"; L="1.048065288040607"; M="1.057171476683828" }
    14 = @{ B="1.02"; C="1.036355896069064"; D="1.043849270886602"; E="1.044646098668355"; F="1.053794969651256"; I="1.035320743173649"; J="1.04292154676227";  K="1.047391206494409"; L="1.048185128423128"; M="1.057300983053314" }
    15 = @{ B="1.02"; C="1.036462948950304"; D="1.043932102351351"; E="1.044739594331464"; F="1.053894243097149"; I="1.03533851059465";  J="1.042991030808273"; K="1.04745440626269";  L="1.048258971151657"; M="1.057380781132962" }
    16 = @{ B="1.02"; C="1.037086254252652"; D="1.04441432709647";  E="1.045284096915101"; F="1.05447236784021";  I="1.035441456178543"; J="1.043395447993409"; K="1.047822117541242"; L="1.048688845112143"; M="1.057845314939698" }
    17 = @{ B="1.02"; C="1.037477422735617"; D="1.044716909627804"; E="1.045625927591253"; F="1.054835282568614"; I="1.035505618218172"; J="1.043649118723253"; K="1.04805265018713";  L="1.048958561062506"; M="1.058136769475464" }
    18 = @{ B="1.02"; C="1.037705651344661"; D="1.044893434974386"; E="1.045825411656418"; F="1.055047062338043"; I="1.035542893857263"; J="1.043797076879283"; K="1.048187071372888"; L="1.049115905789424"; M="1.058306792970869" }
    19 = @{ B="1.02"; C="1.037783482687455"; D="1.044953631296755"; E="1.045893447495361"; F="1.055119290308234"; I="1.035555578600735"; J="1.043847526135636"; K="1.048232897928522"; L="1.049169560313229"; M="1.058364770457327" }
    20 = @{ B="1.02"; C="1.037435447160507"; D="1.044684441852002"; E="1.045589242041043"; F="1.054796335133186"; I="1.035498749646175"; J="1.043621902630111"; K="1.048027920822896"; L="1.048929620606764"; M="1.058105496775853" }
    21 = @{ B="1.02"; C="1.036304734871964"; D="1.043809684307356"; E="1.04460141889887";  F="1.053747528320331"; I="1.035312243131987"; J="1.042888337337843"; K="1.047360998311253"; L="1.048149837350428"; M="1.057262845626181" }
    22 = @{ B="1.02"; C="1.03559441620908";  D="1.043260005274548"; E="1.043981244565766"; F="1.053088991096058"; I="1.035193640065215"; J="1.042427088284076"; K="1.046941282359052"; L="1.047659779174776"; M="1.056733251852372" }
    23 = @{ B="1.02"; C="1.035970911998998"; D="1.043551370174906"; E="1.044309924124005"; F="1.053438009429338"; I="1.035256640560336"; J="1.042671607454195"; K="1.04716381887817";  L="1.047919546879043"; M="1.057013978851475" }
    24 = @{ B="1.02"; C="1.037454413910138"; D="1.044699112537156"; E="1.045605818353145"; F="1.054813933500862"; I="1.035501853718606"; J="1.043634200422551"; K="1.048039095096028"; L="1.048942697478974"; M="1.058119627491718" }
    25 = @{ B="1.02"; C="1.039178547463047"; D="1.046032335727099"; E="1.047113550503223"; F="1.056414446708347"; I="1.035780563514513"; J="1.044751088203617"; K="1.049053051137499"; L="1.050130943898067"; M="1.059403569418626" }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = [double]$cols[$col]
    }
}
